$d = $word.ActiveDocument

# --- Paragraph 1: add extra spacing around the curly-quoted "Review Tips"
#     phrase and a trailing space at the end of the sentence. ---
$rng1 = $d.Content
$ok1 = $rng1.Find.Execute(
    "إذا كنت ترغب في مراجعة أي من النصائح التي تلقيتها سابقًا، ما عليك سوى كتابة MENU والانتقال إلى “Review Tips”",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "إذا كنت ترغب في مراجعة أي من النصائح التي تلقيتها سابقًا، ما عليك سوى كتابة MENU والانتقال إلى   “Review Tips” ",
    2
)
if (-not $ok1) { throw "Replacement 1 (Review Tips paragraph) did not match." }

# --- Paragraph 2: was left in English ("To change your language ...");
#     replace with the (duplicated) Arabic "Review Tips" sentence, matching
#     the translation actually committed upstream. ---
$rng2 = $d.Content
$ok2 = $rng2.Find.Execute(
    "To change your language or gender settings, select “Change my Settings”",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "إذا كنت ترغب في مراجعة أي من النصائح التي تلقيتها سابقًا، ما عليك سوى كتابة MENU والانتقال إلى   “Review Tips”",
    2
)
if (-not $ok2) { throw "Replacement 2 (Change my Settings paragraph) did not match." }

# --- Paragraph 3: translate "To share a link..." sentence to Arabic
#     (keeping the following DOCPROPERTY field run intact). ---
$rng3 = $d.Content
$ok3 = $rng3.Find.Execute(
    "To share a link to this chatbot with a friend, select “Invite a Friend to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "لمشاركة رابط هذا الروبوت مع صديق، اختر “Invite ذ Friend to ",
    2
)
if (-not $ok3) { throw "Replacement 3 (Invite a Friend paragraph) did not match." }

# --- Paragraph 4: translate "For more information..." sentence to Arabic
#     (tail of the sentence, after the closing curly quote, stays English). ---
$rng4 = $d.Content
$ok4 = $rng4.Find.Execute(
    "For more information or resources available to you in a crisis, select “Get more help.” You can also access this information by typing HELP at any time. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "للحصول على مزيد من المعلومات أو الموارد المتاحة لك في أوقات الأزمات، اختر“Get more help.” You can also access this information by typing HELP at any time. ",
    2
)
if (-not $ok4) { throw "Replacement 4 (Get more help paragraph) did not match." }

Write-Output "Replacement 1: $ok1"
Write-Output "Replacement 2: $ok2"
Write-Output "Replacement 3: $ok3"
Write-Output "Replacement 4: $ok4"
